# Refresh the "cryptos" price/volume snapshot on Sheet1 (rows 2-51) to the
# latest scrape: updated Price (D) / Volume(1h) (E) figures for most coins,
# plus a reordering of the NEARProtocol/Filecoin rows (30-31) with their own
# refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks like a plain number need an explicit
# Text number format before assignment, otherwise Excel auto-converts
# the string into a numeric cell (losing the exact textual formatting,
# e.g. trailing zeros or very small decimals rendered in scientific form).
$textCells = @("D5", "D6", "D8", "D11", "D13", "D14", "D18", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '68.463.32'
$ws.Range('E2').Value = '  +2.98%  '
$ws.Range('D3').Value = '3.618.85'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('E4').Value = '  -0.36%  '
$ws.Range('D5').Value = '203.05'
$ws.Range('E5').Value = '  +10.01%  '
$ws.Range('D6').Value = '568.47'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('D7').Value = '3.613.05'
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').Value = '0.624'
$ws.Range('E8').Value = '  +1.57%  '
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('E10').Value = '  +1.34%  '
$ws.Range('D11').Value = '61.19'
$ws.Range('E11').Value = '  +15.96%  '
$ws.Range('E12').Value = '  +4.35%  '
$ws.Range('D13').Value = '0.0000289'
$ws.Range('E13').Value = '  +11.35%  '
$ws.Range('D14').Value = '10.07'
$ws.Range('E14').Value = '  +2.95%  '
$ws.Range('D15').Value = '4.189.21'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '3.603.41'
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('E17').Value = '  +0.98%  '
$ws.Range('D18').Value = '19.08'
$ws.Range('E18').Value = '  +4.22%  '
$ws.Range('D19').Value = '68.224.92'
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').Value = '12.42'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').Value = '404.52'
$ws.Range('D23').Value = '13.40'
$ws.Range('E23').Value = '  +20.31%  '
$ws.Range('E24').Value = '  -3.35%  '
$ws.Range('D25').Value = '85.65'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '2.95'
$ws.Range('E26').Value = '  +1.62%  '
$ws.Range('D27').Value = '12.66'
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('D28').Value = '3.92'
$ws.Range('E28').Value = '  +11.11%  '
$ws.Range('E29').Value = '  +1.42%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '9.40'
$ws.Range('E30').Value = '  +5.33%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = '8.11'
$ws.Range('E31').Value = '  +14.84%  '
$ws.Range('D32').Value = '31.71'
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('D33').Value = '671.22'
$ws.Range('E33').Value = '  +7.88%  '
$ws.Range('D34').Value = '12.32'
$ws.Range('E34').Value = '  +1.23%  '
$ws.Range('E35').Value = '  +1.58%  '
$ws.Range('D36').Value = '64.00'
$ws.Range('E36').Value = '  +0.37%  '
$ws.Range('D37').Value = '42.37'
$ws.Range('E37').Value = '  +2.76%  '
$ws.Range('D38').Value = '0.421'
$ws.Range('E38').Value = '  +5.95%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').Value = '0.0₃0787'
$ws.Range('E40').Value = '  +2.28%  '
$ws.Range('E41').Value = '  +14.17%  '
$ws.Range('D42').Value = '3.230.11'
$ws.Range('E42').Value = '  +8.94%  '
$ws.Range('D43').Value = '0.136'
$ws.Range('E43').Value = '  +3.72%  '
$ws.Range('D44').Value = '2.78'
$ws.Range('E44').Value = '  +11.26%  '
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '2.98'
$ws.Range('E46').Value = '  +25.74%  '
$ws.Range('D47').Value = '2.85'
$ws.Range('E47').Value = '  +14.02%  '
$ws.Range('D48').Value = '0.0421'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('D49').Value = '8.90'
$ws.Range('E49').Value = '  +4.95%  '
$ws.Range('E50').Value = '  +1.06%  '
$ws.Range('D51').Value = '3.07'
$ws.Range('E51').Value = '  -2.50%  '
